$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Remove the first slide (the "Preprocess / Process / PostProcess /
#    Conditions" flow-diagram slide). The two remaining slides shift up,
#    keeping their own content untouched.
# ---------------------------------------------------------------------------
$p.Slides.Item(1).Delete()

# ---------------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" footer date shown on the Slide
#    Master and every Slide Layout (10/15/2019 -> 1/7/2020).
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "1/7/2020"
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes $layouts.Item($i).Shapes
}
